# The deck's slide 16 has a 2-column summary table whose style was changed
# (PowerPoint Table Design gallery -> a different built-in table style).
# This corresponds to the <a:tableStyleId> GUID swap in the OOXML:
#   {4BA9B659-B4B4-4867-8184-CD1EDB0FF732}  ->  {6A69D2C4-D7A9-498E-90CE-3E54125DC14A}
#
# Table styles can't be poked through a plain property set in this object
# model ("Table styles cannot be assigned through a property"), so we go
# through Table.ApplyStyle(StyleId), same as PowerPoint does when a user
# clicks a style swatch.

$p = $ppt.ActivePresentation

# Slide 16, 3rd shape on the slide (title textbox, picture, then the table
# graphicFrame) holds the "Total Outflow" / "All expenses added together"
# table.
$slide = $p.Slides.Item(16)
$tableShape = $slide.Shapes.Item(3)

$table = $tableShape.Table
$table.ApplyStyle("{6A69D2C4-D7A9-498E-90CE-3E54125DC14A}")
